# Scheduled-runner update: refresh cached market-board price snapshots
# (currentAveragePrice* / LevePrice* / LeveProfit*) for a batch of leves
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 486.54544
$ws.Range("I11").Value = 486.54544
$ws.Range("K11").Value = 486.54544
$ws.Range("M11").Value = -346.54544

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2134.75
$ws.Range("J17").Value = 2134.75
$ws.Range("L17").Value = 6404.25
$ws.Range("N17").Value = -6740.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1877
$ws.Range("I40").Value = 1422.25
$ws.Range("K40").Value = 1422.25
$ws.Range("M40").Value = -1247.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 2650.1667
$ws.Range("J42").Value = 5276.3335
$ws.Range("L42").Value = 15829.0005
$ws.Range("N42").Value = -16289.0005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4700
$ws.Range("J64").Value = 4800
$ws.Range("L64").Value = 4800
$ws.Range("N64").Value = -5296

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4700
$ws.Range("J67").Value = 4800
$ws.Range("L67").Value = 4800
$ws.Range("N67").Value = -6516

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 25749.75
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 25749.75
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 25749.75
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -26561.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 25749.75
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 25749.75
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 25749.75
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -28557.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1634
$ws.Range("I98").Value = 887.5714
$ws.Range("K98").Value = 887.5714
$ws.Range("M98").Value = 610.4286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2085.5715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1634
$ws.Range("I122").Value = 887.5714
$ws.Range("K122").Value = 2662.7142
$ws.Range("M122").Value = -212.7142000000003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2466.423
$ws.Range("I138").Value = 1105
$ws.Range("K138").Value = 3315
$ws.Range("M138").Value = 1825

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2038.1818
$ws.Range("I32").Value = 1898.75
$ws.Range("K32").Value = 1898.75
$ws.Range("M32").Value = -1611.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 20000
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4551.2
$ws.Range("J63").Value = 4551.2
$ws.Range("L63").Value = 4551.2
$ws.Range("N63").Value = -5923.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4551.2
$ws.Range("J66").Value = 4551.2
$ws.Range("L66").Value = 22756
$ws.Range("N66").Value = -29620

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 99999
$ws.Range("J53").Value = 99999
$ws.Range("L53").Value = 99999
$ws.Range("N53").Value = -101147

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 48789
$ws.Range("J81").Value = 48789
$ws.Range("L81").Value = 48789
$ws.Range("N81").Value = -50911

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 48789
$ws.Range("J84").Value = 48789
$ws.Range("L84").Value = 146367
$ws.Range("N84").Value = -156975

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 250
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 250
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 100
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2879.6
$ws.Range("I62").Value = 2300
$ws.Range("K62").Value = 2300
$ws.Range("M62").Value = -1676

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2879.6
$ws.Range("I65").Value = 2300
$ws.Range("K65").Value = 11500
$ws.Range("M65").Value = -8380

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2314.8462
$ws.Range("I134").Value = 2372.0908
$ws.Range("K134").Value = 7116.2724
$ws.Range("M134").Value = -4581.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 7875.6
$ws.Range("J39").Value = 9219.5
$ws.Range("L39").Value = 27658.5
$ws.Range("N39").Value = -28246.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1609.1052
$ws.Range("I55").Value = 495
$ws.Range("J55").Value = 1671
$ws.Range("K55").Value = 1485
$ws.Range("L55").Value = 5013
$ws.Range("M55").Value = -1308
$ws.Range("N55").Value = -5367

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 339741
$ws.Range("I128").Value = 339741
$ws.Range("K128").Value = 1019223
$ws.Range("M128").Value = -1014243

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4375.5
$ws.Range("I113").Value = 4375.5
$ws.Range("K113").Value = 4375.5
$ws.Range("M113").Value = -2205.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3027.4167
$ws.Range("I122").Value = 2881.9333
$ws.Range("J122").Value = 3269.889
$ws.Range("K122").Value = 8645.7999
$ws.Range("L122").Value = 9809.667000000001
$ws.Range("M122").Value = -6195.7999
$ws.Range("N122").Value = -14709.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 15554.143
$ws.Range("I68").Value = 1373.3334
$ws.Range("K68").Value = 1373.3334
$ws.Range("M68").Value = -624.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 15554.143
$ws.Range("I71").Value = 1373.3334
$ws.Range("K71").Value = 6866.666999999999
$ws.Range("M71").Value = -3122.666999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10006
$ws.Range("I132").Value = 6008.4
$ws.Range("K132").Value = 18025.2
$ws.Range("M132").Value = -15495.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3075.4167
$ws.Range("I136").Value = 3013.25
$ws.Range("J136").Value = 3199.75
$ws.Range("K136").Value = 9039.75
$ws.Range("L136").Value = 9599.25
$ws.Range("M136").Value = -6489.75
$ws.Range("N136").Value = -14699.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 55000
$ws.Range("J70").Value = 55000
$ws.Range("L70").Value = 55000
$ws.Range("N70").Value = -55630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 55000
$ws.Range("J73").Value = 55000
$ws.Range("L73").Value = 55000
$ws.Range("N73").Value = -57184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
